# Update BOMS: add rows to BOM_CSaAF worksheet describing new parts, and
# switch active sheet/selection to reflect the edits made on that sheet.

$wb = $excel.ActiveWorkbook

$wsCSaAF = $wb.Worksheets.Item("BOM_CSaAF")
$wsSupport = $wb.Worksheets.Item("BOM_NI-USB7856-Support")

# Add the new BOM rows (Qty, Parts, Description) to BOM_CSaAF.
# Values are entered in the same order the original author used (Qty/Parts
# for each row first, then the Description column), so that new shared
# strings land at the same indices as the canonical workbook.
$wsCSaAF.Range("A6").Value = 1
$wsCSaAF.Range("B6").Value = "National Instruments 777584-01 "
$wsCSaAF.Range("C6").Value = "PS-2 Power  Supply for USB7856"

$wsCSaAF.Range("A7").Value = 1
$wsCSaAF.Range("B7").Value = "National Instruments 196739-01"

$wsCSaAF.Range("A8").Value = 1
$wsCSaAF.Range("B8").Value = "National Instruments 196375-01"

$wsCSaAF.Range("C7").Value = "NI 9976 2 Pos Terminal Block for USB7856 PS-2 Power Supply"
$wsCSaAF.Range("C8").Value = "Backshell for NI 9976 Terminal Block"

# Update selection on each sheet
$wsSupport.Range("F26").Select() | Out-Null
$wsCSaAF.Range("C9").Select() | Out-Null

# Make BOM_CSaAF the active (visible/selected) sheet/tab
$wsCSaAF.Activate()
